$d = $word.ActiveDocument

# Locate the "Skillnader" heading paragraph (the anchor point for this edit).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs($i).Range.Text.Trim()
    if ($paraText -eq "Skillnader") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    # Fallback: just use the last paragraph of the document.
    $targetIndex = $d.Paragraphs.Count
}

$anchor = $d.Paragraphs($targetIndex)

# Insert a brand-new paragraph right after the heading.
$anchor.Range.InsertParagraphAfter()

# The freshly created paragraph is now immediately after the anchor;
# give it the new text and the "First Paragraph" style.
$newPara = $d.Paragraphs($targetIndex + 1)
$newPara.Range.Text = "AAa"
$newPara.Style = "First Paragraph"
